$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ScenarioMapping")

# Add the new test case row (TC_011) in row 12
$ws.Range("A12").Value = "TC_011"
$ws.Range("B12").Value = "login.feature"
$ws.Range("C12").Value = "User verify that the logged-in user is admin or support staff"

# D12 uses the same "darker" style as column D above (style used by D3:D11),
# so copy formatting from D11 before setting the value.
$ws.Range("D11").Copy()
$ws.Range("D12").PasteSpecial(-4122)
$ws.Range("D12").Value = "Yes"

$ws.Range("E12").Value = "No"

# Flip SmokeTest column (D3:D11) from Yes to No for the existing rows
$ws.Range("D3:D11").Value = "No"

# Rebuild the data validation rules so the "Yes,No" list dropdown also covers
# the newly inserted D12 cell, and the "Yes" list (E1) stays as-is.
$ws.Cells.Validation.Delete()
$ws.Range("E1").Validation.Add(3, 1, 1, '"Yes"')
$ws.Range("D1:D2").Validation.Add(3, 1, 1, '"Yes,No"')
$ws.Range("D3:D8").Validation.Add(3, 1, 1, '"Yes,No"')
$ws.Range("D9:D14").Validation.Add(3, 1, 1, '"Yes,No"')
$ws.Range("E2:E11").Validation.Add(3, 1, 1, '"Yes,No"')
$ws.Range("E13:E14").Validation.Add(3, 1, 1, '"Yes,No"')
$ws.Range("D15:E21").Validation.Add(3, 1, 1, '"Yes,No"')

# Move the active selection to C15
$ws.Range("C15").Select()
